$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coin table cells to match the refreshed market snapshot.
# D (Price) and E (Volume 1h) columns hold numeric-looking text; force
# Text number format first so Excel COM does not silently coerce them
# into real numbers/percentages (matching the original inline-string cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.90%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.02%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.103"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.30%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.63%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.974"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.68%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.199"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.94%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.931"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.99%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9288"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.94%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1442"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.13%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1953"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.07%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09137"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.24%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03510"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.54%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09841"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.00%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001415"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.73%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005892"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.63%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.594"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.40%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.467"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.14%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1314"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.09%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.802"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.20%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2435"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.47%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04459"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.34%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.56%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004837"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.63%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02100"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05109"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.59%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007471"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.80%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01011"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.85%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.74%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01048"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.88%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006212"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.25%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003059"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-100.00%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001604"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.42%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
